# Update countries & provincias Spain
# - Refresh the covid-19 dashboard counters for a handful of countries
# - "Moldavia" overtook "Afganistan" in total cases, so the two countries
#   swap places in the ranking (rows 65/66)
# - "Islas Malvinas" and "Montserrat" swap places as well (rows 214/215)
# - Bump the "last updated" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 6297847
$ws.Range("C4").Value = 7110
$ws.Range("D4").Value = 3547979
$ws.Range("E4").Value = 2559672
$ws.Range("G4").Value = 232
$ws.Range("H4").Value = 190196

# --- Italia (row 22) ---
$ws.Range("B22").Value = 272912
$ws.Range("C22").Value = 1397
$ws.Range("D22").Value = 208490
$ws.Range("E22").Value = 28915
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = 35507

# --- Alemania (row 23) ---
$ws.Range("B23").Value = 248016
$ws.Range("C23").Value = 625
$ws.Range("E23").Value = 15521
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 9395

# --- Portugal (row 51) ---
$ws.Range("B51").Value = 59051
$ws.Range("C51").Value = 418
$ws.Range("D51").Value = 42427
$ws.Range("E51").Value = 14795
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 1829

# --- Singapur (row 52) ---
$ws.Range("D52").Value = 56028
$ws.Range("E52").Value = 853

# --- Afganistan / Moldavia swap (rows 65/66) ---
# Row 65 keeps the "Afganistan" slot's old numbers in A-column position,
# but it is now Moldavia's row (Moldavia overtook Afganistan).
$ws.Range("A65").Value = "Moldavia"
$ws.Range("B65").Value = 38372
$ws.Range("C65").Value = 632
$ws.Range("D65").Value = 27017
$ws.Range("E65").Value = 10319
$ws.Range("G65").Value = 12
$ws.Range("H65").Value = 1036

$ws.Range("A66").Value = "Afganistan"
$ws.Range("B66").Value = 38288
$ws.Range("C66").Value = 45
$ws.Range("D66").Value = 29390
$ws.Range("E66").Value = 7489
$ws.Range("H66").Value = 1409

# --- Republica de Macedonia (row 86) ---
$ws.Range("D86").Value = 11965
$ws.Range("E86").Value = 2191

# --- Namibia (row 103) ---
$ws.Range("B103").Value = 8082
$ws.Range("C103").Value = 238
$ws.Range("D103").Value = 3483
$ws.Range("E103").Value = 4513
$ws.Range("G103").Value = 4
$ws.Range("H103").Value = 86

# --- Sri Lanka (row 128) ---
$ws.Range("B128").Value = 3111
$ws.Range("C128").Value = 10
$ws.Range("E128").Value = 216

# --- Trinidad yTobago (row 147) ---
$ws.Range("B147").Value = 1941
$ws.Range("C147").Value = 21
$ws.Range("D147").Value = 700
$ws.Range("E147").Value = 1213

# --- Montserrat / Islas Malvinas swap (rows 214/215) ---
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

# --- Timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 17:18"
